$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -------------------------------------------------------------------
# Table 1 (rows 4-7): B = time, C = time (Loop-level parallelism)
# -------------------------------------------------------------------
$ws.Range("B4").Value = 337
$ws.Range("C4").Value = 169

$ws.Range("B5").Value = 255
$ws.Range("C5").Value = 113

$ws.Range("B6").Value = 232
$ws.Range("C6").Value = 123

$ws.Range("B7").Value = 229
$ws.Range("C7").Value = 155

# -------------------------------------------------------------------
# Table 2 (rows 14-17): B = time, C = time (functional decomposition)
# -------------------------------------------------------------------
$ws.Range("B14").Value = 337
$ws.Range("C14").Value = 205

$ws.Range("B15").Value = 255
$ws.Range("C15").Value = 172

$ws.Range("B16").Value = 232
$ws.Range("C16").Value = 178

$ws.Range("B17").Value = 229
$ws.Range("C17").Value = 179

# -------------------------------------------------------------------
# Table 3 (rows 24-27): B = sequential, C = functional, D = loop-level
# -------------------------------------------------------------------
$ws.Range("B24").Value = 337
$ws.Range("C24").Value = 205
$ws.Range("D24").Value = 169

$ws.Range("B25").Value = 255
$ws.Range("C25").Value = 172
$ws.Range("D25").Value = 113

$ws.Range("B26").Value = 232
$ws.Range("C26").Value = 178
$ws.Range("D26").Value = 123

$ws.Range("B27").Value = 229
$ws.Range("C27").Value = 179
$ws.Range("D27").Value = 155

# -------------------------------------------------------------------
# Update the view's selection / scroll position to match the saved file
# -------------------------------------------------------------------
$ws.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 2

$wb.Save()
